$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3.332866452911979
$ws.Range("E2").Value = 16.74055822901289
$ws.Range("F2").Value = 18.3586991835932
$ws.Range("G2").Value = 21.97987378217732
$ws.Range("H2").Value = 10.62994996385521
$ws.Range("O2").Value = 15.51357751095987

# Row 3
$ws.Range("D3").Value = 3.331906116541789
$ws.Range("E3").Value = 15.77358072534923
$ws.Range("F3").Value = 17.85750361967088
$ws.Range("G3").Value = 20.88400046637141
$ws.Range("H3").Value = 10.56398091616973
$ws.Range("O3").Value = 15.21194444154729

# Row 4
$ws.Range("D4").Value = 3.331667210034164
$ws.Range("E4").Value = 15.15375536400613
$ws.Range("F4").Value = 17.54855801098144
$ws.Range("G4").Value = 20.18742174270467
$ws.Range("H4").Value = 10.52596340044021
$ws.Range("O4").Value = 15.02856969616993

# Row 5
$ws.Range("D5").Value = 3.331658326440736
$ws.Range("E5").Value = 14.89487218464951
$ws.Range("F5").Value = 17.42257312030376
$ws.Range("G5").Value = 19.8980169984635
$ws.Range("H5").Value = 10.51110920613299
$ws.Range("O5").Value = 14.95441454222162

# Row 6
$ws.Range("D6").Value = 3.331662202570367
$ws.Range("E6").Value = 14.85151310981643
$ws.Range("F6").Value = 17.40165469098881
$ws.Range("G6").Value = 19.84964027438432
$ws.Range("H6").Value = 10.50868157673871
$ws.Range("O6").Value = 14.9421389230355

# Row 7
$ws.Range("D7").Value = 3.331666731627135
$ws.Range("E7").Value = 15.1502890936753
$ws.Range("F7").Value = 17.54685899359077
$ws.Range("G7").Value = 20.1835405728957
$ws.Range("H7").Value = 10.52576047143518
$ws.Range("O7").Value = 15.02756715204569

# Row 8
$ws.Range("D8").Value = 3.332462701585336
$ws.Range("E8").Value = 16.4127029344815
$ws.Range("F8").Value = 18.18626524397516
$ws.Range("G8").Value = 21.60716863689968
$ws.Range("H8").Value = 10.60669544107147
$ws.Range("O8").Value = 15.40926152103459

# Row 9
$ws.Range("D9").Value = 3.336794265068177
$ws.Range("E9").Value = 18.83526645892443
$ws.Range("F9").Value = 19.42152583551498
$ws.Range("G9").Value = 24.19549305282733
$ws.Range("H9").Value = 10.78454097625014
$ws.Range("O9").Value = 16.1674281179029

# Row 10
$ws.Range("D10").Value = 3.341648870266812
$ws.Range("E10").Value = 20.5158063941146
$ws.Range("F10").Value = 20.30649208104569
$ws.Range("G10").Value = 25.9564039169725
$ws.Range("H10").Value = 10.92598972804871
$ws.Range("O10").Value = 16.72400890826518

# Row 11
$ws.Range("D11").Value = 3.344216520351168
$ws.Range("E11").Value = 21.23749084471423
$ws.Range("F11").Value = 20.70207008789821
$ws.Range("G11").Value = 26.72430725654523
$ws.Range("H11").Value = 10.99247610726912
$ws.Range("O11").Value = 16.97584131324035

# Row 12
$ws.Range("D12").Value = 3.345240104828352
$ws.Range("E12").Value = 21.50465745367271
$ws.Range("F12").Value = 20.85070168403503
$ws.Range("G12").Value = 27.01015662864139
$ws.Range("H12").Value = 11.01794236505452
$ws.Range("O12").Value = 17.07090818264912

# Row 13
$ws.Range("D13").Value = 3.34501738373661
$ws.Range("E13").Value = 21.44738995832759
$ws.Range("F13").Value = 20.81874541673819
$ws.Range("G13").Value = 26.94881594057337
$ws.Range("H13").Value = 11.01244518410014
$ws.Range("O13").Value = 17.05044857741605

# Row 14
$ws.Range("D14").Value = 3.344299706018385
$ws.Range("E14").Value = 21.25959324810562
$ws.Range("F14").Value = 20.71432231239456
$ws.Range("G14").Value = 26.74792412138484
$ws.Range("H14").Value = 10.99456554138557
$ws.Range("O14").Value = 16.9836690268625

# Row 15
$ws.Range("D15").Value = 3.343866772780873
$ws.Range("E15").Value = 21.14376642488968
$ws.Range("F15").Value = 20.65020392082526
$ws.Range("G15").Value = 26.62422428068378
$ws.Range("H15").Value = 10.98365086186755
$ws.Range("O15").Value = 16.94272308504471

# Row 16
$ws.Range("D16").Value = 3.341488264673332
$ws.Range("E16").Value = 20.46778435323411
$ws.Range("F16").Value = 20.28048541154646
$ws.Range("G16").Value = 25.90553718702318
$ws.Range("H16").Value = 10.921686192136
$ws.Range("O16").Value = 16.70751495472316

# Row 17
$ws.Range("D17").Value = 3.340120871657708
$ws.Range("E17").Value = 20.04215383634668
$ws.Range("F17").Value = 20.05176529754217
$ws.Range("G17").Value = 25.45602378938267
$ws.Range("H17").Value = 10.88420815208811
$ws.Range("O17").Value = 16.56280069720596

# Row 18
$ws.Range("D18").Value = 3.339368223788727
$ws.Range("E18").Value = 19.79331164051829
$ws.Range("F18").Value = 19.91956150706801
$ws.Range("G18").Value = 25.19436406922756
$ws.Range("H18").Value = 10.86285440835428
$ws.Range("O18").Value = 16.47944270756333

# Row 19
$ws.Range("D19").Value = 3.339119213877981
$ws.Range("E19").Value = 19.70836458320049
$ws.Range("F19").Value = 19.87469300417123
$ws.Range("G19").Value = 25.10524197005103
$ws.Range("H19").Value = 10.85565975489257
$ws.Range("O19").Value = 16.45120130807224

# Row 20
$ws.Range("D20").Value = 3.340262933356759
$ws.Range("E20").Value = 20.08787984548768
$ws.Range("F20").Value = 20.07618141573259
$ws.Range("G20").Value = 25.50419863862301
$ws.Range("H20").Value = 10.88817691407267
$ws.Range("O20").Value = 16.57821917134397

# Row 21
$ws.Range("D21").Value = 3.344509117089498
$ws.Range("E21").Value = 21.31491949976725
$ws.Range("F21").Value = 20.74502673075179
$ws.Range("G21").Value = 26.80706618550233
$ws.Range("H21").Value = 10.99980952428331
$ws.Range("O21").Value = 17.00329262939371

# Row 22
$ws.Range("D22").Value = 3.347582896093317
$ws.Range("E22").Value = 22.08120714169241
$ws.Range("F22").Value = 21.17530009767941
$ws.Range("G22").Value = 27.62971804028941
$ws.Range("H22").Value = 11.07444526050036
$ws.Range("O22").Value = 17.27933077334006

# Row 23
$ws.Range("D23").Value = 3.345915169710981
$ws.Range("E23").Value = 21.67547574966464
$ws.Range("F23").Value = 20.94633133725009
$ws.Range("G23").Value = 27.19334216298972
$ws.Range("H23").Value = 11.03446353276822
$ws.Range("O23").Value = 17.1321982015476

# Row 24
$ws.Range("D24").Value = 3.340198602956227
$ws.Range("E24").Value = 20.06722001980251
$ws.Range("F24").Value = 20.06514508651023
$ws.Range("G24").Value = 25.48242883184157
$ws.Range("H24").Value = 10.88638203454042
$ws.Range("O24").Value = 16.57124896922865

# Row 25
$ws.Range("D25").Value = 3.335327506432642
$ws.Range("E25").Value = 18.17820302606692
$ws.Range("F25").Value = 19.09055749939828
$ws.Range("G25").Value = 23.51897299350301
$ws.Range("H25").Value = 10.73446794358163
$ws.Range("O25").Value = 15.96195547730651
